$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("1er Parcial")
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws3 = $wb.Worksheets.Item("Final")

# ---------------------------------------------------------------
# "1er Parcial" sheet updates
# ---------------------------------------------------------------
$ws1.Range("C7").Value  = 0
$ws1.Range("D7").Value  = 0
$ws1.Range("E7").Value  = 52
$ws1.Range("F7").Value  = 57.14
$ws1.Range("G7").Value  = 39
$ws1.Range("H7").Value  = 42.86
$ws1.Range("I7").Value  = 7

$ws1.Range("B16").Value = 211
$ws1.Range("C16").Value = 40
$ws1.Range("D16").Value = 18.96
$ws1.Range("E16").Value = 171
$ws1.Range("F16").Value = 81.04000000000001

$ws1.Range("C20").Value = 8
$ws1.Range("D20").Value = 4.3
$ws1.Range("E20").Value = 150
$ws1.Range("F20").Value = 80.65000000000001
$ws1.Range("G20").Value = 28
$ws1.Range("H20").Value = 15.05
$ws1.Range("I20").Value = 7.4

$ws1.Range("C32").Value = 77
$ws1.Range("D32").Value = 31.17
$ws1.Range("E32").Value = 170
$ws1.Range("F32").Value = 68.83
$ws1.Range("I32").Value = 7.8

$ws1.Range("C34").Value = 23
$ws1.Range("D34").Value = 10.36
$ws1.Range("G34").Value = 9
$ws1.Range("H34").Value = 4.05

$ws1.Range("B42").Value = 40
$ws1.Range("C42").Value = 12
$ws1.Range("D42").Value = 30
$ws1.Range("E42").Value = 28
$ws1.Range("F42").Value = 70
$ws1.Range("I42").Value = 8.300000000000001

$ws1.Range("B49").Value = 90
$ws1.Range("C49").Value = 40
$ws1.Range("D49").Value = 44.44
$ws1.Range("E49").Value = 50
$ws1.Range("F49").Value = 55.56
$ws1.Range("I49").Value = 8.6

# ---------------------------------------------------------------
# "2o Parcial" sheet updates
# ---------------------------------------------------------------
$ws2.Range("G7").Value  = 91
$ws2.Range("H7").Value  = 100

$ws2.Range("B16").Value = 211
$ws2.Range("C16").Value = 211
$ws2.Range("G16").Value = 171
$ws2.Range("H16").Value = 81.04000000000001

$ws2.Range("G20").Value = 178
$ws2.Range("H20").Value = 95.7

$ws2.Range("G32").Value = 170
$ws2.Range("H32").Value = 68.83

$ws2.Range("G34").Value = 199
$ws2.Range("H34").Value = 89.64

$ws2.Range("B42").Value = 40
$ws2.Range("C42").Value = 40
$ws2.Range("G42").Value = 28
$ws2.Range("H42").Value = 70

$ws2.Range("B49").Value = 90
$ws2.Range("C49").Value = 90
$ws2.Range("G49").Value = 50
$ws2.Range("H49").Value = 55.56

# ---------------------------------------------------------------
# "Final" sheet updates (mirrors "1er Parcial")
# ---------------------------------------------------------------
$ws3.Range("C7").Value  = 0
$ws3.Range("D7").Value  = 0
$ws3.Range("E7").Value  = 52
$ws3.Range("F7").Value  = 57.14
$ws3.Range("G7").Value  = 39
$ws3.Range("H7").Value  = 42.86
$ws3.Range("I7").Value  = 7

$ws3.Range("B16").Value = 211
$ws3.Range("C16").Value = 40
$ws3.Range("D16").Value = 18.96
$ws3.Range("E16").Value = 171
$ws3.Range("F16").Value = 81.04000000000001

$ws3.Range("C20").Value = 8
$ws3.Range("D20").Value = 4.3
$ws3.Range("E20").Value = 150
$ws3.Range("F20").Value = 80.65000000000001
$ws3.Range("G20").Value = 28
$ws3.Range("H20").Value = 15.05
$ws3.Range("I20").Value = 7.4

$ws3.Range("C32").Value = 77
$ws3.Range("D32").Value = 31.17
$ws3.Range("E32").Value = 170
$ws3.Range("F32").Value = 68.83
$ws3.Range("I32").Value = 7.8

$ws3.Range("C34").Value = 23
$ws3.Range("D34").Value = 10.36
$ws3.Range("G34").Value = 9
$ws3.Range("H34").Value = 4.05

$ws3.Range("B42").Value = 40
$ws3.Range("C42").Value = 12
$ws3.Range("D42").Value = 30
$ws3.Range("E42").Value = 28
$ws3.Range("F42").Value = 70
$ws3.Range("I42").Value = 8.300000000000001

$ws3.Range("B49").Value = 90
$ws3.Range("C49").Value = 40
$ws3.Range("D49").Value = 44.44
$ws3.Range("E49").Value = 50
$ws3.Range("F49").Value = 55.56
$ws3.Range("I49").Value = 8.6
